# Update countries & provincias Spain
# Applies the 7-Jul-2020 20:04 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap the two pairs of countries whose row order changed in the source data ---
# Suazilandia / Mozambique (rows 136-137)
$ws.Cells.Item(136, 1).Value = "Suazilandia"
$ws.Cells.Item(137, 1).Value = "Mozambique"

# Groenlandia / Islas Malvinas (rows 209-210) - identical stats, only label order swaps
$ws.Cells.Item(209, 1).Value = "Groenlandia"
$ws.Cells.Item(210, 1).Value = "Islas Malvinas"

# --- Refresh the "last updated" timestamp ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 7 de Julio de 2020 a las 20:04"

# --- Update the statistics (Casos totales, Nuevos casos, Casos activos, Recuperados,
#     Casos criticos, Muertes hoy, Muertes) for the countries with refreshed figures ---

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

Set-Row 4   3061925 21733 1337373 1591107 0 466 133445   # Estados Unidos
Set-Row 6   742016  21670 456710  264663  0 469 20643    # India
Set-Row 19  198243  186   182700  6449    0 2   9094     # Alemania
Set-Row 23  106106  171   69827   27571   0 15  8708     # Canada
Set-Row 54  25538   7     23364   432     0 1   1742     # Irlanda
Set-Row 66  14607   228   10639   3728    0 3   240      # Marruecos
Set-Row 85  6315    53    4965    1297    0 0   53       # Tayikistan
Set-Row 92  5178    124   2119    3038    0 1   21       # Guayana Francesa
Set-Row 105 3015    9     1096    1827    0 0   92       # Somalia
Set-Row 108 2501    10    2158    331     0 0   12       # Maldivas
Set-Row 110 2395    15    2240    69      0 0   86       # Cuba
Set-Row 128 1297    13    591     358     0 3   348      # Yemen
Set-Row 131 1169    2     969     190     0 0   10       # Jordania
Set-Row 136 1056    45    570     472     0 1   14       # Suazilandia (after label swap)
Set-Row 137 1012    0     277     727     0 0   8        # Mozambique (after label swap)
Set-Row 143 907     66    320     570     0 3   17       # Montenegro
Set-Row 155 551     1     472     76      0 1   3        # Reunion
